# df_RSE_settings.xlsx - "New version of data feeds work but problem with
# the industry balance."
#
# Updates the Minimum/Maximum coefficient pairs for several industries
# (rows 2, 8, 11, 12, 16, 17, 18, 22) and moves the active selection to C14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data updates (column B = Minimum, column C = Maximum) ---
$ws.Range("C2").Value = 5

$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 1

$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 1

$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 1

$ws.Range("B16").Value = 5
$ws.Range("C16").Value = 1

$ws.Range("B17").Value = 5
$ws.Range("C17").Value = 1

$ws.Range("B18").Value = 5
$ws.Range("C18").Value = 1

$ws.Range("B22").Value = 5
$ws.Range("C22").Value = 1

# --- Move the active cell / selection to C14 ---
$ws.Range("C14").Select() | Out-Null
